# Applies the "Globale Planning" update:
#  - Several "Tijd" (time) cells in column D change from "3 uur 15 min" to "3 uur"
#  - The active selection on the Planning sheet moves from C17 to K21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")

$rows = @(37, 57, 59, 64, 66, 69, 70, 74, 75, 77, 78, 82, 84, 86, 87, 94, 95, 98, 100)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "3 uur"
}

$ws.Activate()
$ws.Range("K21").Select()
